$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1 with its text first
$ws.Range("H1").Value = "Save"

# Copy the formatting from the neighboring header cell (G1, which carries the
# shared bold/border/centered header style) onto H1 so it reuses the same
# cell style rather than creating a brand-new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add numeric "Save" values (0) in H2 and H3 beneath the new header
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
